# Generate Report for Handback
#
# The underlying handback-status generator re-ran and produced an updated
# report: the file "05c86cd9-a5e4-4d2a-8084-50674b136977.md" failed its
# handback transform (status flips from "Ready for handoff" to
# "Handback transform failed"), which re-sorts/re-positions it (and the
# now-promoted "4de642fb...md", which was previously "In Translation" and
# swaps into its old slot) in the per-language detail sheets as well as the
# Overview summary sheet. Two other already-"Ready for handoff" rows
# (41fa2a1e...md and f491a28a...md) also swap positions.
#
# Concretely, on every sheet, rows 4 & 5 swap their file/status content, and
# rows 7 & 8 swap their file content (status unchanged, both stay "Ready for
# handoff").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: columns A (File Name), B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Row 4 <-> Row 5 (file name + status both change: 05c86cd9 now failed,
# 4de642fb now takes the "In Translation" slot)
$overview.Range("A4").Value = "05c86cd9-a5e4-4d2a-8084-50674b136977.md"
$overview.Range("B4").Value = "Handback transform failed"
$overview.Range("C4").Value = "Handback transform failed"

$overview.Range("A5").Value = "4de642fb-5d28-458a-b184-8a8e132fd194.md"
$overview.Range("B5").Value = "In Translation"
$overview.Range("C5").Value = "In Translation"

# Row 7 <-> Row 8 (status stays "Ready for handoff" on both; only the file
# name moves)
$overview.Range("A7").Value = "41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"
$overview.Range("A8").Value = "f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"

# ---------------------------------------------------------------------------
# Per-language detail sheets: zh-cn & de-de
# Columns: A Source File Name, B Status, C Latest Handoff File,
#          D Latest Handoff Datetime, ... H Handoff Reason
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Row4Handoff = "05c86cd9-a5e4-4d2a-8084-50674b136977.111de57fd2c0d2429ff87b39452ba18ea590b9aa.zh-cn.xlf"; Row4Date = "2016-03-03 16:01:45"; Row5Handoff = "4de642fb-5d28-458a-b184-8a8e132fd194.87ff7d79b1e73c865135837a6d7e438a54468aa6.zh-cn.xlf"; Row5Date = "2016-03-03 15:58:09"; Row7Handoff = "41fa2a1e-e5f6-419a-8cae-3684c8394aac.e5e53929549bef51e53423cb6de6c93f89a0399f.zh-cn.xlf"; Row8Handoff = "f491a28a-ae0e-4d0d-98aa-0ad501f29e48.1755ebbbfef550e4347980f9ae77e572d2349b51.zh-cn.xlf" },
    @{ Name = "de-de"; Row4Handoff = "05c86cd9-a5e4-4d2a-8084-50674b136977.111de57fd2c0d2429ff87b39452ba18ea590b9aa.de-de.xlf"; Row4Date = "2016-03-03 16:02:11"; Row5Handoff = "4de642fb-5d28-458a-b184-8a8e132fd194.87ff7d79b1e73c865135837a6d7e438a54468aa6.de-de.xlf"; Row5Date = "2016-03-03 15:58:23"; Row7Handoff = "41fa2a1e-e5f6-419a-8cae-3684c8394aac.e5e53929549bef51e53423cb6de6c93f89a0399f.de-de.xlf"; Row8Handoff = "f491a28a-ae0e-4d0d-98aa-0ad501f29e48.1755ebbbfef550e4347980f9ae77e572d2349b51.de-de.xlf" }
)

foreach ($info in $langSheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Row 4 becomes the failed-handback file (05c86cd9)
    $ws.Range("A4").Value = "05c86cd9-a5e4-4d2a-8084-50674b136977.md"
    $ws.Range("B4").Value = "Handback transform failed"
    $ws.Range("C4").Value = $info.Row4Handoff
    $ws.Range("D4").Value = $info.Row4Date

    # Row 5 becomes the in-translation file (4de642fb)
    $ws.Range("A5").Value = "4de642fb-5d28-458a-b184-8a8e132fd194.md"
    $ws.Range("B5").Value = "In Translation"
    $ws.Range("C5").Value = $info.Row5Handoff
    $ws.Range("D5").Value = $info.Row5Date

    # Row 7 becomes 41fa2a1e (still "Ready for handoff")
    $ws.Range("A7").Value = "41fa2a1e-e5f6-419a-8cae-3684c8394aac.md"
    $ws.Range("C7").Value = $info.Row7Handoff

    # Row 8 becomes f491a28a (still "Ready for handoff")
    $ws.Range("A8").Value = "f491a28a-ae0e-4d0d-98aa-0ad501f29e48.md"
    $ws.Range("C8").Value = $info.Row8Handoff
}
